$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCols = @("J", "K", "L", "M", "N", "O")

# ---------------------------------------------------------------------------
# Row 1: new header cells J1:O1, styled like the existing header cells.
# ---------------------------------------------------------------------------
$headerValues = @(
    "URL",
    "Aktivitetsnummer",
    "Rapporttittel",
    "Dato",
    "Oppgaveleder",
    "Deltakere_i_revisjon"
)

for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value = $headerValues[$i]
}

# Copy just the formatting (bold/centered/bordered) from the existing header
# cell I1 onto each new header cell, so they share the same style index
# instead of minting new, duplicate ones.
$ws.Range("I1").Copy()
foreach ($col in $newCols) {
    $ws.Range($col + "1").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Rows 3-11: new, empty placeholder cells J:O - matching the pre-existing
# empty placeholder cells already present in columns F:I on those rows.
# ---------------------------------------------------------------------------
for ($r = 3; $r -le 11; $r++) {
    foreach ($col in $newCols) {
        $addr = $col + $r
        $ws.Range($addr).Value = "'"
        $ws.Range("F" + $r).Copy()
        $ws.Range($addr).PasteSpecial(-4122)
    }
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 2: new data values J2:O2.
# ---------------------------------------------------------------------------
$ws.Range("J2").Value = "https://www.ptil.no//contentassets/f955a9317bff4749bfc21db88383a8a3/2020_102_rapport-palegg-tilsyn-neptune-gjoa-vedlikeholdsstyring.pdf"

# K2 looks numeric ("027153044 ") but must stay text, preserving the leading
# zero and the trailing space - a leading apostrophe forces text entry.
$ws.Range("K2").Value = "'027153044 "

$ws.Range("L2").Value = "undervannsanleggene (aktivitet 027153044)  "
$ws.Range("M2").Value = "27.04.2020 "
$ws.Range("N2").Value = "Kenneth Skogen "
$ws.Range("O2").Value = "Mihajlovic "

# Re-normalize K2's style back to the plain (un-styled) look of its row-mates
# (J2, L2, M2, N2, O2), now that the text value has been safely written.
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
